# Update onboarding survey with alternate formats
# - replaces the nested "Please rate the following statements" likert group
#   (select_one kt2eu04, appearance list-nolabel) with a flat set of
#   acknowledge/note/select_one fields, each using its own choice list and
#   its own "alternate format" appearance (horizontal-compact, likert,
#   compact, quickcompact).
# - adds the corresponding choice lists (ml6mb78, om8ix13, ak2us99, vj3vs28,
#   wq0ck21) to the choices sheet.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------
# survey sheet: rows 18-26 (old) become rows 18-25 (new)
# ---------------------------------------------------------------------

# Row 18 used to open a nested "field-list" group; now it just closes the
# group_mf7bp00 group that was opened on row 16.
$survey.Range("A18:G18").ClearContents()
$survey.Cells.Item(18, 1).Value = "end_group"

# Row 19 used to be the select_one kt2eu04 header row; now it's a plain
# acknowledge field.
$survey.Range("A19:G19").ClearContents()
$survey.Cells.Item(19, 1).Value = "acknowledge"
$survey.Cells.Item(19, 2).Value = "Acknowledge"
$survey.Cells.Item(19, 3).Value = "Acknowledge"
$survey.Cells.Item(19, 5).Value = "false"

# Insert a new row 20 for the "note" field introduced before the statements.
$survey.Rows.Item(20).Insert()
$survey.Cells.Item(20, 1).Value = "note"
$survey.Cells.Item(20, 2).Value = "For_each_of_the_foll_e_with_the_statement"
$survey.Cells.Item(20, 3).Value = "For each of the following statements, please indicate the extent to which you agree with the statement."
$survey.Cells.Item(20, 5).Value = "false"

# Row 21 (was row 20): first likert statement -> select_one ml6mb78
$survey.Range("A21:G21").ClearContents()
$survey.Cells.Item(21, 1).Value = "select_one ml6mb78"
$survey.Cells.Item(21, 2).Value = "The_information_you_l_your_driving_needs"
$survey.Cells.Item(21, 3).Value = "The information you have about electric vehicles and charging electric infrastructure is enough to fulfill your driving needs"
$survey.Cells.Item(21, 5).Value = "false"

# Row 22 (was row 21): select_one om8ix13, appearance horizontal-compact,
# slightly reworded label.
$survey.Range("A22:G22").ClearContents()
$survey.Cells.Item(22, 1).Value = "select_one om8ix13"
$survey.Cells.Item(22, 2).Value = "You_are_concerned_th_charging_opportunity"
$survey.Cells.Item(22, 3).Value = "You are concerned that an electric vehicle could run out of charge before the next charging opportunity"
$survey.Cells.Item(22, 5).Value = "false"
$survey.Cells.Item(22, 7).Value = "horizontal-compact"

# Row 23 (was row 22): select_one ak2us99, appearance likert.
$survey.Range("A23:G23").ClearContents()
$survey.Cells.Item(23, 1).Value = "select_one ak2us99"
$survey.Cells.Item(23, 2).Value = "Concern_of_running_o_ectric_fleet_vehicle"
$survey.Cells.Item(23, 3).Value = "Concern of running out of charge, or range anxiety, prevents you from choosing an electric fleet vehicle"
$survey.Cells.Item(23, 5).Value = "false"
$survey.Cells.Item(23, 7).Value = "likert"

# Row 24 (was row 23): select_one vj3vs28, appearance compact.
$survey.Range("A24:G24").ClearContents()
$survey.Cells.Item(24, 1).Value = "select_one vj3vs28"
$survey.Cells.Item(24, 2).Value = "You_believe_that_the_an_electric_vehicle"
$survey.Cells.Item(24, 3).Value = "You believe that the number of electric vehicles charging stations are few within your area of operation and it is a significant obstacle to driving an electric vehicle"
$survey.Cells.Item(24, 5).Value = "false"
$survey.Cells.Item(24, 7).Value = "compact"

# Row 25 (was row 24): select_one wq0ck21, appearance quickcompact.
$survey.Range("A25:G25").ClearContents()
$survey.Cells.Item(25, 1).Value = "select_one wq0ck21"
$survey.Cells.Item(25, 2).Value = "The_use_of_public_ra_nimize_range_anxiety"
$survey.Cells.Item(25, 3).Value = "The use of public rapid-charging infrastructure as a back-up option would minimize range anxiety"
$survey.Cells.Item(25, 5).Value = "false"
$survey.Cells.Item(25, 7).Value = "quickcompact"

# Old rows 25 & 26 held the two closing end_group rows (closing the nested
# "Please_rate_the_following_statements" group and group_mf7bp00 itself).
# group_mf7bp00 is now closed directly on row 18, and the statements are no
# longer wrapped in a group, so both trailing end_group rows (now at 26 and
# 27 after the insert on row 20) are obsolete and must be removed.
$survey.Rows.Item("26:27").Delete()

# ---------------------------------------------------------------------
# choices sheet: replace the kt2eu04 list (rows 29-33) with ml6mb78, and
# append the om8ix13 / ak2us99 / vj3vs28 / wq0ck21 copies (rows 34-53).
# ---------------------------------------------------------------------

$likertNames = @("disagree", "somewhat_disagree", "neither_agree_nor_disagree", "somewhat_agree", "agree")
$likertLabels = @("Disagree", "Somewhat disagree", "Neither agree nor disagree", "Somewhat agree", "Agree")
$likertLists = @("ml6mb78", "om8ix13", "ak2us99", "vj3vs28", "wq0ck21")

$row = 29
foreach ($listName in $likertLists) {
    for ($i = 0; $i -lt 5; $i++) {
        $choices.Cells.Item($row, 1).Value = $listName
        $choices.Cells.Item($row, 2).Value = $likertNames[$i]
        $choices.Cells.Item($row, 3).Value = $likertLabels[$i]
        $row = $row + 1
    }
}
